# Client: Fix Weapon Delete Bug
# Update the "Move" sheet stamina-cost values for rows 7 and 8,
# and move the active selection to C15 as left by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Move")
$ws.Activate()

$ws.Range("C7").Value = 0.8
$ws.Range("C8").Value = 1

$ws.Range("C15").Select()
